$d = $word.ActiveDocument

# Append a new paragraph "This is nonsence." at the very end of the
# document. The first letter needs its own run carrying an
# eastAsia-hinted <w:rFonts/>, matching the same split the preceding
# paragraph already uses ("I" / " am adding..."), so we build the
# paragraph from explicit OOXML rather than relying on typed-text runs.
$xml = '<?xml version="1.0"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:r>' +
                '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>' +
                '<w:t>T</w:t>' +
              '</w:r>' +
              '<w:r>' +
                '<w:t>his is nonsence.</w:t>' +
              '</w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$target = $d.Content
$target.Collapse(0)
$target.InsertXML($xml)
